$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.684.54'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.586.19'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('E4').Value = '  +1.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.38'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('E7').Value = '  +1.32%  '
$ws.Range('E8').Value = '  -3.70%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.811.54'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.588.78'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.89'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('E15').Value = '  -4.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.650.22'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.33'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.63'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0695'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('E20').Value = '  -4.05%  '
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('E22').Value = '  -3.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.57'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.53%  '
$ws.Range('E24').Value = '  -2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.86'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.12'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('E29').Value = '  -3.56%  '
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0467'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.380.81'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.61%  '
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.975'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.537'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.978'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.37%  '
$ws.Range('E43').Value = '  -2.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.54'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.721.87'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.30'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('E49').Value = '  +11.77%  '
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('E51').Value = '  -0.59%  '
